$wb = $excel.ActiveWorkbook

# --- HW Architecture sheet ---
$wsHw = $wb.Worksheets.Item("HW Architecture")
$wsHw.Activate()
$wsHw.Range("D2:D5").Value = "12.5G"
$wsHw.Range("D6:D69").Value = "1.4G"
$wsHw.Range("D6:D69").Select()

# --- Availability sheet ---
$wsAvail = $wb.Worksheets.Item("Availability")
$wsAvail.Activate()
$wsAvail.Range("D3").Value = 24
$wsAvail.Range("D4").Value = 24
$wsAvail.Range("D5").Value = 24
$wsAvail.Range("D6").Value = 24
$wsAvail.Range("D8").Value = 24
$wsAvail.Range("D9").Value = 24
$wsAvail.Range("H27").Select()

# --- Manual Calculation sheet (becomes the active tab) ---
$wsCalc = $wb.Worksheets.Item("Manual Calculation")
$wsCalc.Activate()
$wsCalc.Range("I36").Select()
